$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 338.14285
$ws.Range("I33").Value = 165.77777
$ws.Range("J33").Value = 648.4
$ws.Range("K33").Value = 165.77777
$ws.Range("L33").Value = 648.4
$ws.Range("M33").Value = 63.22223
$ws.Range("N33").Value = -1106.4
$ws.Range("H138").Value = 2337.0286
$ws.Range("I138").Value = 1434.2858
$ws.Range("J138").Value = 2938.8572
$ws.Range("K138").Value = 4302.857400000001
$ws.Range("L138").Value = 8816.571599999999
$ws.Range("M138").Value = 837.1425999999992
$ws.Range("N138").Value = -19096.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30326
$ws.Range("J44").Value = 30326
$ws.Range("L44").Value = 30326
$ws.Range("N44").Value = -31302
$ws.Range("H55").Value = 23602
$ws.Range("J55").Value = 23602
$ws.Range("L55").Value = 23602
$ws.Range("N55").Value = -24232
$ws.Range("H80").Value = 34075
$ws.Range("J80").Value = 34075
$ws.Range("L80").Value = 34075
$ws.Range("N80").Value = -36071
$ws.Range("H83").Value = 34075
$ws.Range("J83").Value = 34075
$ws.Range("L83").Value = 102225
$ws.Range("N83").Value = -112209
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29285.5
$ws.Range("J35").Value = 31326.285
$ws.Range("L35").Value = 31326.285
$ws.Range("N35").Value = -31946.285
$ws.Range("H82").Value = 15962.6
$ws.Range("I82").Value = 3071.4285
$ws.Range("J82").Value = 27242.375
$ws.Range("K82").Value = 3071.4285
$ws.Range("L82").Value = 27242.375
$ws.Range("M82").Value = -2688.4285
$ws.Range("N82").Value = -28008.375
$ws.Range("H85").Value = 15962.6
$ws.Range("I85").Value = 3071.4285
$ws.Range("J85").Value = 27242.375
$ws.Range("K85").Value = 3071.4285
$ws.Range("L85").Value = 27242.375
$ws.Range("M85").Value = -1745.4285
$ws.Range("N85").Value = -29894.375
$ws.Range("H99").Value = 1958.0385
$ws.Range("I99").Value = 1274.9375
$ws.Range("J99").Value = 3051
$ws.Range("K99").Value = 1274.9375
$ws.Range("L99").Value = 3051
$ws.Range("M99").Value = 223.0625
$ws.Range("N99").Value = -6047

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19232754
$ws.Range("I31").Value = 47620496
$ws.Range("J31").Value = 2348.6128
$ws.Range("K31").Value = 47620496
$ws.Range("L31").Value = 2348.6128
$ws.Range("M31").Value = -47620201
$ws.Range("N31").Value = -2938.6128
$ws.Range("H34").Value = 19232754
$ws.Range("I34").Value = 47620496
$ws.Range("J34").Value = 2348.6128
$ws.Range("K34").Value = 47620496
$ws.Range("L34").Value = 2348.6128
$ws.Range("M34").Value = -47620294
$ws.Range("N34").Value = -2752.6128
$ws.Range("H41").Value = 22012.6
$ws.Range("J41").Value = 22012.6
$ws.Range("L41").Value = 22012.6
$ws.Range("N41").Value = -22868.6
$ws.Range("H51").Value = 29555
$ws.Range("J51").Value = 29555
$ws.Range("L51").Value = 29555
$ws.Range("N51").Value = -31027
$ws.Range("H60").Value = 12898.6
$ws.Range("J60").Value = 17650
$ws.Range("L60").Value = 17650
$ws.Range("N60").Value = -18672
$ws.Range("H61").Value = 29555
$ws.Range("J61").Value = 29555
$ws.Range("L61").Value = 29555
$ws.Range("N61").Value = -30251
$ws.Range("H68").Value = 25867.857
$ws.Range("J68").Value = 25867.857
$ws.Range("L68").Value = 25867.857
$ws.Range("N68").Value = -27365.857
$ws.Range("H71").Value = 25867.857
$ws.Range("J71").Value = 25867.857
$ws.Range("L71").Value = 77603.571
$ws.Range("N71").Value = -85091.571
$ws.Range("H109").Value = 21286.6
$ws.Range("J109").Value = 21286.6
$ws.Range("L109").Value = 21286.6
$ws.Range("N109").Value = -23366.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 923.05554
$ws.Range("I131").Value = 366.66666
$ws.Range("J131").Value = 1034.3334
$ws.Range("K131").Value = 1099.99998
$ws.Range("L131").Value = 3103.0002
$ws.Range("M131").Value = 3940.00002
$ws.Range("N131").Value = -13183.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 21459.385
$ws.Range("J57").Value = 22997.666
$ws.Range("L57").Value = 22997.666
$ws.Range("N57").Value = -24637.666
$ws.Range("H122").Value = 2852.8462
$ws.Range("I122").Value = 2546.8333
$ws.Range("J122").Value = 3115.1428
$ws.Range("K122").Value = 7640.499899999999
$ws.Range("L122").Value = 9345.428400000001
$ws.Range("M122").Value = -5190.499899999999
$ws.Range("N122").Value = -14245.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 142858780
$ws.Range("I16").Value = 250001380
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 250001380
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -250001210
$ws.Range("N16").Value = -2340
$ws.Range("H76").Value = 15144
$ws.Range("J76").Value = 15144
$ws.Range("L76").Value = 15144
$ws.Range("N76").Value = -15820
$ws.Range("H79").Value = 15144
$ws.Range("J79").Value = 15144
$ws.Range("L79").Value = 15144
$ws.Range("N79").Value = -17484
$ws.Range("H122").Value = 3484.9333
$ws.Range("I122").Value = 2944.3333
$ws.Range("J122").Value = 4295.8335
$ws.Range("K122").Value = 8832.999899999999
$ws.Range("L122").Value = 12887.5005
$ws.Range("M122").Value = -6382.999899999999
$ws.Range("N122").Value = -17787.5005
$ws.Range("H132").Value = 3667.1667
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 4250.75
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 12752.25
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -17812.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 31030.8
$ws.Range("J109").Value = 31030.8
$ws.Range("L109").Value = 31030.8
$ws.Range("N109").Value = -33804.8
$ws.Range("H132").Value = 2543.8076
$ws.Range("I132").Value = 2650.6667
$ws.Range("J132").Value = 2223.2307
$ws.Range("K132").Value = 7952.000100000001
$ws.Range("L132").Value = 6669.6921
$ws.Range("M132").Value = -5422.000100000001
$ws.Range("N132").Value = -11729.6921
